$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-06-13 Thursday" "2024-06-14 Friday"
Replace-Text "801÷4=" "509÷7="
Replace-Text "195÷3=" "643÷3="
Replace-Text "770÷3=" "772÷3="
Replace-Text "156÷3=" "523÷6="
Replace-Text "240÷6=" "223÷8="
Replace-Text "186÷7=" "364÷6="
Replace-Text "796÷7=" "526÷7="
Replace-Text "812÷8=" "792÷7="
Replace-Text "228÷4=" "660÷9="
Replace-Text "757÷5=" "406÷6="
Replace-Text "360÷4=" "665÷3="
Replace-Text "177÷6=" "727÷4="
Replace-Text "872÷4=" "757÷4="
Replace-Text "735÷9=" "118÷6="
Replace-Text "650÷9=" "129÷4="
Replace-Text "835÷5=" "471÷8="
Replace-Text "403÷8=" "421÷9="
Replace-Text "959÷2=" "631÷7="
Replace-Text "612÷4=" "646÷4="
Replace-Text "709÷2=" "367÷4="
Replace-Text "541÷8=" "275÷9="
Replace-Text "661÷8=" "758÷3="
Replace-Text "319÷5=" "445÷2="
Replace-Text "782÷2=" "905÷8="
Replace-Text "139÷5=" "318÷5="
